# Apply the changes described by the commit:
#  - Remove the stray "unclear" values duplicated in column B for rows 18 & 19
#    (column C already holds the correct "unclear" value for those rows).
#  - Move the active cell selection from C1 to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicated column B entries on rows 18 and 19.
$ws.Range("B18").ClearContents()
$ws.Range("B19").ClearContents()

# Update the saved selection to B1.
$ws.Range("B1").Select()
